# Auto-generated: update cached market-price / profit values per scheduled-runner diff.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 4554.2856
$ws.Range("I38").Value = 460
$ws.Range("J38").Value = 5670.909
$ws.Range("K38").Value = 1380
$ws.Range("L38").Value = 17012.727
$ws.Range("M38").Value = -1008
$ws.Range("N38").Value = -17756.727
$ws.Range("H40").Value = 3094.9473
$ws.Range("I40").Value = 1000
$ws.Range("J40").Value = 3653.6
$ws.Range("K40").Value = 1000
$ws.Range("L40").Value = 3653.6
$ws.Range("M40").Value = -825
$ws.Range("N40").Value = -4003.6
$ws.Range("H43").Value = 2289.48
$ws.Range("I43").Value = 603.4375
$ws.Range("J43").Value = 5286.8887
$ws.Range("K43").Value = 603.4375
$ws.Range("L43").Value = 5286.8887
$ws.Range("M43").Value = -534.4375
$ws.Range("N43").Value = -5424.8887
$ws.Range("H98").Value = 6164.081
$ws.Range("I98").Value = 4282.0713
$ws.Range("K98").Value = 4282.0713
$ws.Range("M98").Value = -2784.0713
$ws.Range("H122").Value = 6164.081
$ws.Range("I122").Value = 4282.0713
$ws.Range("K122").Value = 12846.2139
$ws.Range("M122").Value = -10396.2139
$ws.Range("H123").Value = 41787.145
$ws.Range("J123").Value = 41787.145
$ws.Range("L123").Value = 41787.145
$ws.Range("N123").Value = -51587.145
$ws.Range("H129").Value = 1355.8551
$ws.Range("J129").Value = 1372.6119
$ws.Range("L129").Value = 4117.835700000001
$ws.Range("N129").Value = -14117.8357
$ws.Range("H137").Value = 520086.97
$ws.Range("I137").Value = 1255638.6
$ws.Range("J137").Value = 2476.537
$ws.Range("K137").Value = 3766915.8
$ws.Range("L137").Value = 7429.610999999999
$ws.Range("M137").Value = -3764365.8
$ws.Range("N137").Value = -12529.611
$ws.Range("H138").Value = 3202.612
$ws.Range("I138").Value = 1673.8
$ws.Range("J138").Value = 3853.1702
$ws.Range("K138").Value = 5021.4
$ws.Range("L138").Value = 11559.5106
$ws.Range("M138").Value = 118.6000000000004
$ws.Range("N138").Value = -21839.5106

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3386.862
$ws.Range("I32").Value = 3206.164
$ws.Range("J32").Value = 3992.2
$ws.Range("K32").Value = 3206.164
$ws.Range("L32").Value = 3992.2
$ws.Range("M32").Value = -2919.164
$ws.Range("N32").Value = -4566.2
$ws.Range("H45").Value = 1586.375
$ws.Range("I45").Value = 1679.4166
$ws.Range("K45").Value = 1679.4166
$ws.Range("M45").Value = -1302.4166
$ws.Range("H122").Value = 3460.6487
$ws.Range("I122").Value = 3129.138
$ws.Range("K122").Value = 9387.414000000001
$ws.Range("M122").Value = -6937.414000000001
$ws.Range("H132").Value = 1820.2727
$ws.Range("I132").Value = 1023.9643
$ws.Range("J132").Value = 6279.6
$ws.Range("K132").Value = 3071.8929
$ws.Range("L132").Value = 18838.8
$ws.Range("M132").Value = -541.8928999999998
$ws.Range("N132").Value = -23898.8
$ws.Range("H139").Value = 42642.5
$ws.Range("J139").Value = 42642.5
$ws.Range("L139").Value = 42642.5
$ws.Range("N139").Value = -52922.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H133").Value = 39772.668
$ws.Range("J133").Value = 57900
$ws.Range("L133").Value = 57900
$ws.Range("N133").Value = -68020
$ws.Range("H138").Value = 40736.43
$ws.Range("J138").Value = 40736.43
$ws.Range("L138").Value = 40736.43
$ws.Range("N138").Value = -51016.43

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2043.197
$ws.Range("I31").Value = 876.4483
$ws.Range("J31").Value = 2957.6758
$ws.Range("K31").Value = 876.4483
$ws.Range("L31").Value = 2957.6758
$ws.Range("M31").Value = -581.4483
$ws.Range("N31").Value = -3547.6758
$ws.Range("H34").Value = 2043.197
$ws.Range("I34").Value = 876.4483
$ws.Range("J34").Value = 2957.6758
$ws.Range("K34").Value = 876.4483
$ws.Range("L34").Value = 2957.6758
$ws.Range("M34").Value = -674.4483
$ws.Range("N34").Value = -3361.6758
$ws.Range("H58").Value = 2573.9736
$ws.Range("I58").Value = 1518.9286
$ws.Range("J58").Value = 5528.1
$ws.Range("K58").Value = 1518.9286
$ws.Range("L58").Value = 5528.1
$ws.Range("M58").Value = -1315.9286
$ws.Range("N58").Value = -5934.1
$ws.Range("H99").Value = 14289543
$ws.Range("I99").Value = 28573086
$ws.Range("J99").Value = 6000
$ws.Range("K99").Value = 28573086
$ws.Range("L99").Value = 6000
$ws.Range("M99").Value = -28571588
$ws.Range("N99").Value = -8996
$ws.Range("H105").Value = 2029.2858
$ws.Range("I105").Value = 1985
$ws.Range("J105").Value = 2062.5
$ws.Range("K105").Value = 1985
$ws.Range("L105").Value = 2062.5
$ws.Range("M105").Value = -238
$ws.Range("N105").Value = -5556.5
$ws.Range("H126").Value = 14289543
$ws.Range("I126").Value = 28573086
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 85719258
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = -85716788
$ws.Range("N126").Value = -22940
$ws.Range("H134").Value = 1618.3043
$ws.Range("I134").Value = 1010.5238
$ws.Range("J134").Value = 8000
$ws.Range("K134").Value = 3031.5714
$ws.Range("L134").Value = 24000
$ws.Range("M134").Value = -496.5714000000003
$ws.Range("N134").Value = -29070
$ws.Range("H136").Value = 2573.9736
$ws.Range("I136").Value = 1518.9286
$ws.Range("J136").Value = 5528.1
$ws.Range("K136").Value = 4556.7858
$ws.Range("L136").Value = 16584.3
$ws.Range("M136").Value = -2006.7858
$ws.Range("N136").Value = -21684.3

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3413.8462
$ws.Range("I3").Value = 3002
$ws.Range("J3").Value = 4786.6665
$ws.Range("K3").Value = 9006
$ws.Range("L3").Value = 14359.9995
$ws.Range("M3").Value = -8894
$ws.Range("N3").Value = -14583.9995
$ws.Range("H68").Value = 1518.9474
$ws.Range("J68").Value = 1648.7142
$ws.Range("L68").Value = 4946.142599999999
$ws.Range("N68").Value = -6568.142599999999
$ws.Range("H71").Value = 1518.9474
$ws.Range("J71").Value = 1648.7142
$ws.Range("L71").Value = 14838.4278
$ws.Range("N71").Value = -22950.4278
$ws.Range("H98").Value = 439.8
$ws.Range("I98").Value = 999
$ws.Range("J98").Value = 300
$ws.Range("K98").Value = 2997
$ws.Range("L98").Value = 900
$ws.Range("M98").Value = -1499
$ws.Range("N98").Value = -3896
$ws.Range("H122").Value = 2321.4363
$ws.Range("I122").Value = 721.7
$ws.Range("J122").Value = 3235.5715
$ws.Range("K122").Value = 6495.3
$ws.Range("L122").Value = 29120.1435
$ws.Range("M122").Value = -4045.3
$ws.Range("N122").Value = -34020.1435

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 29999
$ws.Range("J51").Value = 29999
$ws.Range("L51").Value = 29999
$ws.Range("N51").Value = -31017
$ws.Range("H102").Value = 3521.261
$ws.Range("I102").Value = 2416.1667
$ws.Range("K102").Value = 2416.1667
$ws.Range("M102").Value = -794.1667000000002
$ws.Range("H126").Value = 3510.6
$ws.Range("I126").Value = 2872.1177
$ws.Range("J126").Value = 4867.375
$ws.Range("K126").Value = 8616.3531
$ws.Range("L126").Value = 14602.125
$ws.Range("M126").Value = -6146.3531
$ws.Range("N126").Value = -19542.125
$ws.Range("H132").Value = 3832
$ws.Range("I132").Value = 2223
$ws.Range("J132").Value = 5226.467
$ws.Range("K132").Value = 6669
$ws.Range("L132").Value = 15679.401
$ws.Range("M132").Value = -4139
$ws.Range("N132").Value = -20739.401

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 20918292
$ws.Range("I22").Value = 25100950
$ws.Range("J22").Value = 5000
$ws.Range("K22").Value = 25100950
$ws.Range("L22").Value = 5000
$ws.Range("M22").Value = -25100655
$ws.Range("N22").Value = -5590
$ws.Range("H27").Value = 20918292
$ws.Range("I27").Value = 25100950
$ws.Range("J27").Value = 5000
$ws.Range("K27").Value = 25100950
$ws.Range("L27").Value = 5000
$ws.Range("M27").Value = -25100843
$ws.Range("N27").Value = -5214
$ws.Range("H46").Value = 1117.3448
$ws.Range("I46").Value = 774.06665
$ws.Range("J46").Value = 1485.1428
$ws.Range("K46").Value = 774.06665
$ws.Range("L46").Value = 1485.1428
$ws.Range("M46").Value = -586.06665
$ws.Range("N46").Value = -1861.1428
$ws.Range("H110").Value = 30333.334
$ws.Range("J110").Value = 30333.334
$ws.Range("L110").Value = 30333.334
$ws.Range("N110").Value = -38513.334
$ws.Range("H122").Value = 5425.1763
$ws.Range("I122").Value = 5077.3335
$ws.Range("J122").Value = 6260
$ws.Range("K122").Value = 15232.0005
$ws.Range("L122").Value = 18780
$ws.Range("M122").Value = -12782.0005
$ws.Range("N122").Value = -23680

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 65185.89
$ws.Range("J46").Value = 65185.89
$ws.Range("L46").Value = 65185.89
$ws.Range("N46").Value = -65647.89
$ws.Range("H133").Value = 35463.637
$ws.Range("J133").Value = 35463.637
$ws.Range("L133").Value = 35463.637
$ws.Range("N133").Value = -45583.637
$ws.Range("H134").Value = 65185.89
$ws.Range("J134").Value = 65185.89
$ws.Range("L134").Value = 195557.67
$ws.Range("N134").Value = -200627.67
$ws.Range("H136").Value = 6061.5625
$ws.Range("I136").Value = 3279
$ws.Range("J136").Value = 7326.364
$ws.Range("K136").Value = 9837
$ws.Range("L136").Value = 21979.092
$ws.Range("M136").Value = -7287
$ws.Range("N136").Value = -27079.092
